# Add reset user data feature
#
# For six Saturdays/working days in the "Annual Leave" column (G) that were
# incorrectly marked, move the 1.0 mark to the "At Work" column (C) instead:
# rows 32, 33, 34, 37, 38, 41 (22/23/24/27/28/31 January 2025).
# Then refresh the Total row (44) so C44 / G44 reflect the six moved marks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(32, 33, 34, 37, 38, 41)

foreach ($r in $rows) {
    # "At Work" (C) gets a 1.0 mark - copy it (value+style) from a row that
    # already carries the same literal "1.0" text/style (row 31), so the
    # inserted text is indistinguishable from the pre-existing marks.
    $ws.Range("C31").Copy($ws.Range("C$r"))

    # "Annual Leave" (G) no longer applies for this day - clear it back to
    # the same blank state as the other working days (keeps style s=10).
    $ws.Range("G$r").ClearContents()
}

# Recompute the Total row: 6 days moved from Annual Leave -> At Work.
# At Work: 14.0 -> 20.0 ; Annual Leave: 6.0 -> "-" (0 -> dash placeholder).
$scratch = $ws.Range("Z100")
$scratch.Formula = "=""20.0"""
$scratch.Copy()
$ws.Range("C44").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.ClearContents()

# G44's "-" placeholder already exists (F44) with the identical Total-row
# style - copy it across instead of typing a fresh literal.
$ws.Range("F44").Copy($ws.Range("G44"))
